$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.803.54'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '2.462.06'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").Value = '2.461.49'
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("D16").Value = '2.907.14'
$ws.Range("D17").Value = '62.708.17'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '2.462.30'
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +10.17%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +20.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '659.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.28%  '
$ws.Range("D28").Value = '0.0₃0979'
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -14.97%  '
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +2.64%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("E39").Value = '  -2.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '151.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").Value = '0.0₆0310'
$ws.Range("E44").Value = '  -62.88%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("E51").Value = '  -1.12%  '
